$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

# Locate the "filename" column from the header row (row 1).
$filenameCol = 0
for ($c = 1; $c -le $colCount; $c++) {
    $header = $ws.Cells.Item(1, $c).Value()
    if ($header -ne $null -and $header.ToString() -eq "filename") {
        $filenameCol = $c
        break
    }
}

if ($filenameCol -eq 0) {
    $filenameCol = 15  # fallback: known column "O"
}

for ($row = 2; $row -le $rowCount; $row++) {
    $cell = $ws.Cells.Item($row, $filenameCol)
    $current = $cell.Value()
    if ($current -ne $null) {
        $text = $current.ToString()
        if ($text.StartsWith("/img/")) {
            $cell.Value = $text.Substring(1)
        }
    }
}
